$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 50000
$ws.Range("I21").Value = 50000
$ws.Range("K21").Value = 50000
$ws.Range("M21").Value = -49532
$ws.Range("H23").Value = 50000
$ws.Range("I23").Value = 50000
$ws.Range("K23").Value = 50000
$ws.Range("M23").Value = -49766
$ws.Range("H38").Value = 298.5
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").Value = $null
$ws.Range("H40").Value = 5968.75
$ws.Range("J40").Value = 7000
$ws.Range("L40").Value = 7000
$ws.Range("N40").Value = -7350
$ws.Range("H101").Value = 497.7
$ws.Range("I101").Value = 476.5
$ws.Range("K101").Value = 1429.5
$ws.Range("M101").Value = 192.5
$ws.Range("H115").Value = 105.333336
$ws.Range("I115").Value = 105.333336
$ws.Range("K115").Value = 316.000008
$ws.Range("M115").Value = 1250.999992

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 484.6
$ws.Range("I2").Value = 484.6
$ws.Range("K2").Value = 484.6
$ws.Range("M2").Value = -371.6
$ws.Range("H5").Value = 291.5
$ws.Range("I5").Value = 291.5
$ws.Range("K5").Value = 291.5
$ws.Range("M5").Value = -179.5
$ws.Range("H32").Value = 1899.7333
$ws.Range("I32").Value = 1620.4138
$ws.Range("K32").Value = 1620.4138
$ws.Range("M32").Value = -1333.4138
$ws.Range("H96").Value = 10000
$ws.Range("I96").Value = 2000
$ws.Range("J96").Value = 18000
$ws.Range("K96").Value = 2000
$ws.Range("L96").Value = 18000
$ws.Range("M96").Value = 746
$ws.Range("N96").Value = -23492
$ws.Range("H98").Value = 19400
$ws.Range("J98").Value = 19400
$ws.Range("L98").Value = 19400
$ws.Range("N98").Value = -25390
$ws.Range("H116").Value = 484.6
$ws.Range("I116").Value = 484.6
$ws.Range("K116").Value = 484.6
$ws.Range("M116").Value = 1809.4
$ws.Range("H121").Value = 40255
$ws.Range("J121").Value = 40255
$ws.Range("L121").Value = 40255
$ws.Range("N121").Value = -43749
$ws.Range("H122").Value = 2544.6667
$ws.Range("J122").Value = 3251.6667
$ws.Range("L122").Value = 9755.000100000001
$ws.Range("N122").Value = -14655.0001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 484.6
$ws.Range("I3").Value = 484.6
$ws.Range("K3").Value = 484.6
$ws.Range("M3").Value = -370.6
$ws.Range("H4").Value = 291.5
$ws.Range("I4").Value = 291.5
$ws.Range("K4").Value = 291.5
$ws.Range("M4").Value = -176.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 11198.667
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 11198.667
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 11198.667
$ws.Range("M17").Value = $null
$ws.Range("N17").Value = -11546.667
$ws.Range("H22").Value = 661.9
$ws.Range("I22").Value = 353.16666
$ws.Range("K22").Value = 353.16666
$ws.Range("M22").Value = -3.166659999999979
$ws.Range("H31").Value = 1531.3704
$ws.Range("I31").Value = 1439.7693
$ws.Range("K31").Value = 1439.7693
$ws.Range("M31").Value = -1144.7693
$ws.Range("H33").Value = 5080
$ws.Range("I33").Value = 2133.3333
$ws.Range("K33").Value = 2133.3333
$ws.Range("M33").Value = -1754.3333
$ws.Range("H34").Value = 1531.3704
$ws.Range("I34").Value = 1439.7693
$ws.Range("K34").Value = 1439.7693
$ws.Range("M34").Value = -1237.7693
$ws.Range("H59").Value = 29389.223
$ws.Range("I59").Value = 19752
$ws.Range("J59").Value = 32142.715
$ws.Range("K59").Value = 19752
$ws.Range("L59").Value = 32142.715
$ws.Range("M59").Value = -18607
$ws.Range("N59").Value = -34432.715
$ws.Range("H60").Value = 21454.555
$ws.Range("I60").Value = 22181.834
$ws.Range("J60").Value = 20000
$ws.Range("K60").Value = 22181.834
$ws.Range("L60").Value = 20000
$ws.Range("M60").Value = -21670.834
$ws.Range("N60").Value = -21022
$ws.Range("H99").Value = 4731.5713
$ws.Range("I99").Value = 3799.5
$ws.Range("J99").Value = 5104.4
$ws.Range("K99").Value = 3799.5
$ws.Range("L99").Value = 5104.4
$ws.Range("M99").Value = -2301.5
$ws.Range("N99").Value = -8100.4
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").Value = $null
$ws.Range("H126").Value = 4731.5713
$ws.Range("I126").Value = 3799.5
$ws.Range("J126").Value = 5104.4
$ws.Range("K126").Value = 11398.5
$ws.Range("L126").Value = 15313.2
$ws.Range("M126").Value = -8928.5
$ws.Range("N126").Value = -20253.2

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 200
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = $null
$ws.Range("H50").Value = 432.33334
$ws.Range("I50").Value = 432.33334
$ws.Range("K50").Value = 1297.00002
$ws.Range("M50").Value = -816.0000199999999
$ws.Range("H53").Value = 432.33334
$ws.Range("I53").Value = 432.33334
$ws.Range("K53").Value = 1297.00002
$ws.Range("M53").Value = -816.0000199999999
$ws.Range("H55").Value = 1659
$ws.Range("I55").Value = 595
$ws.Range("J55").Value = 1925
$ws.Range("K55").Value = 1785
$ws.Range("L55").Value = 5775
$ws.Range("M55").Value = -1608
$ws.Range("N55").Value = -6129
$ws.Range("H132").Value = 5000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 45000
$ws.Range("M132").Value = $null
$ws.Range("N132").Value = -50060

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 176.45454
$ws.Range("I107").Value = 202.25
$ws.Range("K107").Value = 202.25
$ws.Range("M107").Value = 1717.75
$ws.Range("H132").Value = 2141.4285
$ws.Range("I132").Value = 2038.2
$ws.Range("J132").Value = 2399.5
$ws.Range("K132").Value = 6114.6
$ws.Range("L132").Value = 7198.5
$ws.Range("M132").Value = -3584.6
$ws.Range("N132").Value = -12258.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2328.1428
$ws.Range("I61").Value = 2328.1428
$ws.Range("K61").Value = 2328.1428
$ws.Range("M61").Value = -2126.1428
$ws.Range("H113").Value = 2328.1428
$ws.Range("I113").Value = 2328.1428
$ws.Range("K113").Value = 2328.1428
$ws.Range("M113").Value = -158.1428000000001
$ws.Range("H122").Value = 3444
$ws.Range("I122").Value = 3444
$ws.Range("K122").Value = 10332
$ws.Range("M122").Value = -7882
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").Value = $null
